# Applies the "Erweiterung der Outputdatei: Status" edit to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename header text for column E from "Ampel/Kreuzung?" to "Ampel?"
$ws.Range("E1").Value = "Ampel?"

# 2. Mark the "Bushaltestelle?" status flag (column D) as active for rows 3-10
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1

# 3. Update the target speed for row 6
$ws.Range("B6").Value = 20

# 4. Update the saved selection to match the author's last cursor position
$ws.Range("G9").Select()
